# Dailyscrum update: Kevin (column D) and Bram (column F) fill in their
# "what did I do / what will I do / what problems do I have" entries for
# the 04-03-2016, 11-03-2016 and 18-03-2016 sheets.

$wb = $excel.ActiveWorkbook

$sheet304 = $wb.Worksheets.Item("02-03-2016")
$sheet0403 = $wb.Worksheets.Item("04-03-2016")
$sheet1103 = $wb.Worksheets.Item("11-03-2016")
$sheet1803 = $wb.Worksheets.Item("18-03-2016")

# --- Kevin Strijbos (column D) works through his updates sheet by sheet ---
$sheet0403.Range("D3").Value = "contextdiagram maken, agendapunten voorbereiden en werkjes controleren"
$sheet0403.Range("D4").Value = "statistiekenpagina en homepage mockup"
$sheet0403.Range("D5").Value = "/"

$sheet1103.Range("D3").Value = "statistiekenpagina en homepage mockup"
$sheet1103.Range("D4").Value = "UML class diagram"

$sheet1803.Range("D3").Value = "UML class diagram"

# --- Bram Van Vleymen (column F) works through his updates sheet by sheet ---
$sheet0403.Range("F4").Value = "templates"
$sheet0403.Range("F5").Value = "/"

$sheet1103.Range("F3").Value = "templates"
$sheet1103.Range("F4").Value = "mockup evenementen lijst"

$sheet1803.Range("F3").Value = "mockup evenementen lijst"

# --- leave the cursor/selection where it ended up on each visited sheet ---
$sheet304.Range("D4").Select()
$sheet0403.Range("F5").Select()
$sheet1103.Range("F4").Select()
$sheet1803.Range("F4").Select()
